$d = $word.ActiveDocument

# Update the placeholder ID text in the first paragraph.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5304_topic_10__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SMC_PGI_5304_102__ID**", 2)

# Remove the now-orphaned trailing space that used to live in its own run.
$p = $d.Paragraphs(1)
$rng = $p.Range
$spaceRange = $d.Range($rng.End - 2, $rng.End - 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

# Give the first paragraph a thin paragraph border (space-only, no line) and
# bump its left indent from 120 to 225 twips (6pt -> 11.25pt).
$p = $d.Paragraphs(1)
$p.LeftIndent = 11.25
$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
